$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.488079071044922
$ws.Range("B1").Value = 3.173280954360962
$ws.Range("C1").Value = 5.258664131164551
$ws.Range("D1").Value = 1.550381898880005
$ws.Range("E1").Value = 0.8201225399971008
